$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 317.05554
$ws.Range("I53").Value = 89
$ws.Range("J53").Value = 499.5
$ws.Range("K53").Value = 89
$ws.Range("L53").Value = 499.5
$ws.Range("M53").Value = 548
$ws.Range("N53").Value = -1773.5

$ws.Range("H62").Value = 2700
$ws.Range("I62").Value = 2700
$ws.Range("K62").Value = 2700
$ws.Range("M62").Value = -2076

$ws.Range("H65").Value = 2700
$ws.Range("I65").Value = 2700
$ws.Range("K65").Value = 13500
$ws.Range("M65").Value = -10380

$ws.Range("H132").Value = 2563.3027
$ws.Range("I132").Value = 2297.3174
$ws.Range("J132").Value = 3852.3076
$ws.Range("K132").Value = 6891.9522
$ws.Range("L132").Value = 11556.9228
$ws.Range("M132").Value = -4361.9522
$ws.Range("N132").Value = -16616.9228

$ws.Range("H134").Value = 40714.285
$ws.Range("J134").Value = 40714.285
$ws.Range("L134").Value = 40714.285
$ws.Range("N134").Value = -50854.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2084.07
$ws.Range("I32").Value = 1993.8866
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 1993.8866
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -1706.8866
$ws.Range("N32").Value = -5574

$ws.Range("H122").Value = 2004.1
$ws.Range("I122").Value = 1567.4286
$ws.Range("K122").Value = 4702.2858
$ws.Range("M122").Value = -2252.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8021.4863
$ws.Range("I86").Value = 9447.414000000001
$ws.Range("J86").Value = 2852.5
$ws.Range("K86").Value = 9447.414000000001
$ws.Range("L86").Value = 2852.5
$ws.Range("M86").Value = -8324.414000000001
$ws.Range("N86").Value = -5098.5

$ws.Range("H89").Value = 8021.4863
$ws.Range("I89").Value = 9447.414000000001
$ws.Range("J89").Value = 2852.5
$ws.Range("K89").Value = 47237.07000000001
$ws.Range("L89").Value = 14262.5
$ws.Range("M89").Value = -41621.07000000001
$ws.Range("N89").Value = -25494.5

$ws.Range("H105").Value = 1800.375
$ws.Range("I105").Value = 1730.3846
$ws.Range("J105").Value = 2103.6667
$ws.Range("K105").Value = 1730.3846
$ws.Range("L105").Value = 2103.6667
$ws.Range("M105").Value = 16.61539999999991
$ws.Range("N105").Value = -5597.6667

$ws.Range("H134").Value = 2225.9546
$ws.Range("I134").Value = 2165.3333
$ws.Range("J134").Value = 2698.8
$ws.Range("K134").Value = 6495.999899999999
$ws.Range("L134").Value = 8096.400000000001
$ws.Range("M134").Value = -3960.999899999999
$ws.Range("N134").Value = -13166.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41768.37
$ws.Range("I31").Value = 38778.535
$ws.Range("J31").Value = 43370.07
$ws.Range("K31").Value = 38778.535
$ws.Range("L31").Value = 43370.07
$ws.Range("M31").Value = -38483.535
$ws.Range("N31").Value = -43960.07

$ws.Range("H34").Value = 41768.37
$ws.Range("I34").Value = 38778.535
$ws.Range("J34").Value = 43370.07
$ws.Range("K34").Value = 38778.535
$ws.Range("L34").Value = 43370.07
$ws.Range("M34").Value = -38576.535
$ws.Range("N34").Value = -43774.07

$ws.Range("H58").Value = 21740540
$ws.Range("I58").Value = 31251348
$ws.Range("J58").Value = 1551.6428
$ws.Range("K58").Value = 31251348
$ws.Range("L58").Value = 1551.6428
$ws.Range("M58").Value = -31251145
$ws.Range("N58").Value = -1957.6428

$ws.Range("H99").Value = 1726.1111
$ws.Range("I99").Value = 1691.875
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1691.875
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -193.875
$ws.Range("N99").Value = -4996

$ws.Range("H126").Value = 1726.1111
$ws.Range("I126").Value = 1691.875
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5075.625
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -2605.625
$ws.Range("N126").Value = -10940

$ws.Range("H136").Value = 21740540
$ws.Range("I136").Value = 31251348
$ws.Range("J136").Value = 1551.6428
$ws.Range("K136").Value = 93754044
$ws.Range("L136").Value = 4654.928400000001
$ws.Range("M136").Value = -93751494
$ws.Range("N136").Value = -9754.928400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 209.61539
$ws.Range("I50").Value = 96.875
$ws.Range("J50").Value = 390
$ws.Range("K50").Value = 290.625
$ws.Range("L50").Value = 1170
$ws.Range("M50").Value = 190.375
$ws.Range("N50").Value = -2132

$ws.Range("H53").Value = 209.61539
$ws.Range("I53").Value = 96.875
$ws.Range("J53").Value = 390
$ws.Range("K53").Value = 290.625
$ws.Range("L53").Value = 1170
$ws.Range("M53").Value = 190.375
$ws.Range("N53").Value = -2132

$ws.Range("H75").Value = 1959.5
$ws.Range("I75").Value = 1555
$ws.Range("J75").Value = 2364
$ws.Range("K75").Value = 4665
$ws.Range("L75").Value = 7092
$ws.Range("M75").Value = -3667
$ws.Range("N75").Value = -9088

$ws.Range("H78").Value = 1959.5
$ws.Range("I78").Value = 1555
$ws.Range("J78").Value = 2364
$ws.Range("K78").Value = 13995
$ws.Range("L78").Value = 21276
$ws.Range("M78").Value = -9003
$ws.Range("N78").Value = -31260

$ws.Range("H113").Value = 564.3
$ws.Range("I113").Value = 505.5
$ws.Range("J113").Value = 589.5
$ws.Range("K113").Value = 1516.5
$ws.Range("L113").Value = 1768.5
$ws.Range("M113").Value = 653.5
$ws.Range("N113").Value = -6108.5

$ws.Range("H131").Value = 949.8684
$ws.Range("J131").Value = 1019.24243
$ws.Range("L131").Value = 3057.72729
$ws.Range("N131").Value = -13137.72729

$ws.Range("H132").Value = 2296.5833
$ws.Range("I132").Value = 1313
$ws.Range("J132").Value = 4263.75
$ws.Range("K132").Value = 11817
$ws.Range("L132").Value = 38373.75
$ws.Range("M132").Value = -9287
$ws.Range("N132").Value = -43433.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1109.5135
$ws.Range("I97").Value = 1131.5667
$ws.Range("J97").Value = 1015
$ws.Range("K97").Value = 1131.5667
$ws.Range("L97").Value = 1015
$ws.Range("M97").Value = -635.5667000000001
$ws.Range("N97").Value = -2007

$ws.Range("H102").Value = 1435.6666
$ws.Range("I102").Value = 1425.3334
$ws.Range("J102").Value = 1466.6666
$ws.Range("K102").Value = 1425.3334
$ws.Range("L102").Value = 1466.6666
$ws.Range("M102").Value = 196.6666
$ws.Range("N102").Value = -4710.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 754.8946999999999
$ws.Range("I16").Value = 549.58826
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 549.58826
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -379.58826
$ws.Range("N16").Value = -2840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1939.0769
$ws.Range("I122").Value = 945.3333
$ws.Range("J122").Value = 4175
$ws.Range("K122").Value = 2835.9999
$ws.Range("L122").Value = 12525
$ws.Range("M122").Value = -385.9998999999998
$ws.Range("N122").Value = -17425

$ws.Range("H132").Value = 48261.832
$ws.Range("I132").Value = 37316.445
$ws.Range("J132").Value = 67963.53
$ws.Range("K132").Value = 111949.335
$ws.Range("L132").Value = 203890.59
$ws.Range("M132").Value = -109419.335
$ws.Range("N132").Value = -208950.59
